$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert 3 new columns at F (nav_title, nav_icon_name, nav_icon_lib) ---
$ws.Columns("F:H").Insert()

# --- 2) Insert 1 new column after sect_title_bgcolor (now col T) for sect_title_color ---
$ws.Columns("U").Insert()
# restore the raw column width that was on the neighbour before the insert (~15.664)
$ws.Columns("U").ColumnWidth = 14.83

# --- header row (row 1) new headers ---
$ws.Range("F1").Value = "nav_title"
$ws.Range("G1").Value = "nav_icon_name"
$ws.Range("H1").Value = "nav_icon_lib"
$ws.Range("U1").Value = "sect_title_color"

# --- new nav_title / nav_icon_name values on rows 4 and 5 ---
$ws.Range("F4").Value = "titre1"
$ws.Range("G4").Value = "user"
$ws.Range("F5").Value = "titre2"
$ws.Range("G5").Value = "user"

# --- sect_title_align column (now R) value updates ---
$ws.Range("R2").Value = "start"
$ws.Range("R3").Value = "end"
$ws.Range("R4").Value = "center"
$ws.Range("R5").Value = "center"

# --- clear now-removed sect_title_bgcolor value on row3 (now T3) ---
$ws.Range("T3").Clear()

# --- clear sect_footer_bgcolor column (now V) entirely ---
$ws.Range("V2").Clear()

# --- sect_tlp_msg column (now W): row2 message changes, others cleared ---
$ws.Range("W2").Value = "yeaaaahhhhhh !"
$ws.Range("W3").Clear()
$ws.Range("W4").Clear()
$ws.Range("W5").Clear()
$ws.Range("W6").Clear()
$ws.Range("W7").Clear()
$ws.Range("W8").Clear()

# --- sect_tlp_color column (now X): cleared on every data row ---
$ws.Range("X2").Clear()
$ws.Range("X3").Clear()
$ws.Range("X4").Clear()
$ws.Range("X5").Clear()
$ws.Range("X6").Clear()
$ws.Range("X7").Clear()
$ws.Range("X8").Clear()

# --- the column-insert left behind empty styled placeholder cells where no
#     data exists in the target; remove them so those rows have no <c> there ---
$ws.Range("F2").Clear()
$ws.Range("F3").Clear()
$ws.Range("F6").Clear()
$ws.Range("F7").Clear()
$ws.Range("F8").Clear()
$ws.Range("G2").Clear()
$ws.Range("G3").Clear()
$ws.Range("G6").Clear()
$ws.Range("G7").Clear()
$ws.Range("G8").Clear()
$ws.Range("H2:H8").Clear()
$ws.Range("U2").Clear()
$ws.Range("U3").Clear()

# --- row heights ---
$ws.Range("A1:Y1").RowHeight = 28.8
$ws.Range("A2:Y2").RowHeight = 72

# --- sheet view: scroll / selection ---
$ws.Range("V6").Select()
$excel.ActiveWindow.ScrollColumn = 13

Write-Host "edit applied"
